# ---------------------------------------------------------------------------
# Lab 8 "Entrega Final" edit: rewrite two answer paragraphs in the
# "Preguntas de analisis" section of Docs/Observaciones-Lab 8.docx.
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Question 1 answer: append the missing closing remark about balanced
#    trees, replacing the old closing sentence about "size".
# ---------------------------------------------------------------------------
$old6 = "Están ligados en el sentido de que cada nivel de la altura del árbol contiene la totalidad de los elementos del mismo (size)."
$new6 = "Si el árbol está balanceado, su altura es aproximadamente Log(N) donde N es su tamaño."
$ok6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
if (-not $ok6) { throw "Find/Replace #1 (balanced tree sentence) did not match" }

# ---------------------------------------------------------------------------
# 2) Question 3 answer: the TAD operation used to go from a date range to the
#    list of values is actually "om.values(...)" (Ordered Map), and the
#    final count comes from "lt.size()" (Lista), not "getValue()"/"size()".
#    Replace the whole explanation (keeping the untouched lead-in run and the
#    single leading space before it) and then re-apply bold to the operation
#    names, matching the emphasis used for the other TAD calls in this doc.
# ---------------------------------------------------------------------------
$old12 = " getValue(). Sin embargo, para obtener la información final deseada (número de crímenes entre las fechas) se utiliza la operación size(), dado que esta retorna los elementos que se encuentran en esa lista."
$new12 = " om.values(ordered map, llave_min, llave_max) (del TAD Ordered Map). Esta función recibe un mapa ordenado, y dos llaves: primero la menor de las dos y luego la mayor. La función entonces retorna los valores tales que su llave está contenida en el intervalo inclusivo (o “cerrado”) definido por las dos fechas. Sin embargo, luego en esta funcionalidad del programa para obtener la información final deseada (número de crímenes entre las fechas) se utiliza la operación lt.size() (del TAD Lista), dado que esta retorna los elementos que se encuentran en esa lista."
$rng = $d.Content
$ok12 = $rng.Find.Execute($old12, $true, $false, $false, $false, $false, $true, 1, $false, $new12, 2)
if (-not $ok12) { throw "Find/Replace #2 (om.values/lt.size rewrite) did not match" }

$s = $rng.Start

$bold1 = $d.Range($s + 1, $s + 45)
if ($bold1.Text -ne "om.values(ordered map, llave_min, llave_max)") { throw "bold span #1 text mismatch: [$($bold1.Text)]" }
$bold1.Font.Bold = 1
$bold1.Font.BoldBi = 1

$bold2 = $d.Range($s + 468, $s + 478)
if ($bold2.Text -ne "lt.size() ") { throw "bold span #2 text mismatch: [$($bold2.Text)]" }
$bold2.Font.Bold = 1
$bold2.Font.BoldBi = 1

Write-Output "Lab 8 edits applied."
